# Append two new data rows (197, 198) to the sheet, matching the data
# produced by the author's R script run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DateCell($row, $col, $value) {
    # Copy number format / style from the cell directly above (which
    # already carries the "yyyy-mm-dd hh:mm:ss" date style used by the
    # rest of column A), then overwrite the value.
    $src = $ws.Cells.Item($row - 1, $col)
    $dst = $ws.Cells.Item($row, $col)
    $src.Copy($dst)
    $dst.Value = $value
}

function Set-NumberCell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-TextCell($row, $col, $text) {
    # Force a genuine text cell (t="s") even when the text looks like a
    # number (e.g. "2.87"), then strip the temporary text number-format
    # back off so no stray cell style lingers behind.
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

# ---- Row 197 ----
Set-DateCell   197 1 45481.2916666667
Set-NumberCell 197 2 0
Set-NumberCell 197 3 2.86999988555908
Set-NumberCell 197 4 2.86999988555908
Set-NumberCell 197 5 2.86999988555908
Set-NumberCell 197 6 2.86999988555908
Set-TextCell   197 7 "2.86999988555908"
Set-TextCell   197 8 "EAV.MI"

# ---- Row 198 ----
Set-DateCell   198 1 45482.649375
Set-NumberCell 198 2 18500
Set-NumberCell 198 3 2.79999995231628
Set-NumberCell 198 4 2.57999992370605
Set-NumberCell 198 5 2.79999995231628
Set-NumberCell 198 6 2.67000007629395
Set-TextCell   198 7 "2.67000007629395"
Set-TextCell   198 8 "EAV.MI"

Write-Output "rows 197-198 added"
